# ---------------------------------------------------------------------------
# Applies the "ooutput update 2025 august" commit:
#   1. Metadata sheet: update the canonical URL and the generation Date.
#   2. Elements sheet: shrink a batch of "best fit" column widths (the
#      IG-publisher re-ran with a slightly narrower font/measurement pass),
#      while leaving the untouched columns exactly as they were.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Elements"

# --- 1. Metadata text updates -----------------------------------------------
$ws1.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mindfulness-session-duration"
$ws1.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- 2. Elements sheet column widths ----------------------------------------
# Columns that keep the same width in the diff are left completely untouched
# so their stored width is not perturbed. Columns whose "best fit" width
# shrank are resized via ColumnWidth (character units); the handful that are
# hidden keep their Hidden state re-applied since assigning ColumnWidth
# resets it.

$col1 = $ws2.Columns.Item(1)
$col1.ColumnWidth = 15.714285714285714

$col2 = $ws2.Columns.Item(2)
$col2.ColumnWidth = 15.714285714285714

$col3 = $ws2.Columns.Item(3)
$col3.ColumnWidth = 9.0
$col3.Hidden = $true

$col4 = $ws2.Columns.Item(4)
$col4.ColumnWidth = 6.142857142857143
$col4.Hidden = $true

$col5 = $ws2.Columns.Item(5)
$col5.ColumnWidth = 4.428571428571429

$col6 = $ws2.Columns.Item(6)
$col6.ColumnWidth = 3.142857142857143

$col7 = $ws2.Columns.Item(7)
$col7.ColumnWidth = 3.4285714285714284

$col8 = $ws2.Columns.Item(8)
$col8.ColumnWidth = 11.857142857142858

$col9 = $ws2.Columns.Item(9)
$col9.ColumnWidth = 9.714285714285714

$col11 = $ws2.Columns.Item(11)
$col11.ColumnWidth = 7.428571428571429

$col15 = $ws2.Columns.Item(15)
$col15.ColumnWidth = 11.428571428571429

$col20 = $ws2.Columns.Item(20)
$col20.ColumnWidth = 7.0

$col21 = $ws2.Columns.Item(21)
$col21.ColumnWidth = 12.857142857142858

$col22 = $ws2.Columns.Item(22)
$col22.ColumnWidth = 13.142857142857142

$col23 = $ws2.Columns.Item(23)
$col23.ColumnWidth = 14.142857142857142

$col24 = $ws2.Columns.Item(24)
$col24.ColumnWidth = 13.857142857142858

$col25 = $ws2.Columns.Item(25)
$col25.ColumnWidth = 16.142857142857142

$col26 = $ws2.Columns.Item(26)
$col26.ColumnWidth = 14.285714285714286

$col27 = $ws2.Columns.Item(27)
$col27.ColumnWidth = 4.142857142857143

$col28 = $ws2.Columns.Item(28)
$col28.ColumnWidth = 17.142857142857142

$col29 = $ws2.Columns.Item(29)
$col29.ColumnWidth = 33.714285714285715

$col30 = $ws2.Columns.Item(30)
$col30.ColumnWidth = 12.714285714285714

$col31 = $ws2.Columns.Item(31)
$col31.ColumnWidth = 10.428571428571429
$col31.Hidden = $true

$col32 = $ws2.Columns.Item(32)
$col32.ColumnWidth = 14.142857142857142
$col32.Hidden = $true

$col33 = $ws2.Columns.Item(33)
$col33.ColumnWidth = 7.285714285714286
$col33.Hidden = $true

$col34 = $ws2.Columns.Item(34)
$col34.ColumnWidth = 7.714285714285714

$col37 = $ws2.Columns.Item(37)
$col37.ColumnWidth = 18.714285714285715
